$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the "time" timestamp (column E) for the existing rows 2-31 ---
$timeUpdates = @{
    2 = 44814.77682810185
    3 = 44814.76734414352
    4 = 44814.7760153588
    5 = 44814.77075734954
    6 = 44814.7702228125
    7 = 44816.00728251157
    8 = 44816.00867809028
    9 = 44814.76977856481
    10 = 44814.77070542824
    11 = 44814.7746128125
    12 = 44816.00769685186
    13 = 44814.77207561343
    14 = 44814.77248436343
    15 = 44814.77724549769
    16 = 44814.77501091435
    17 = 44816.00839300926
    18 = 44816.00879736111
    19 = 44814.77160474537
    20 = 44814.77616778935
    21 = 44816.00888778935
    22 = 44814.77166428241
    23 = 44814.77226959491
    24 = 44814.77006715278
    25 = 44816.00448315972
    26 = 44814.77435706019
    27 = 44814.96243638889
    28 = 44814.77139508102
    29 = 44814.77274806713
    30 = 44814.76999071759
    31 = 44814.76836144676
}
foreach ($r in $timeUpdates.Keys) {
    $ws.Range("E$r").Value = $timeUpdates[$r]
}

# --- Append the new sentence rows 32-71 ---
$newRows = @(
    @{ Row=32; A=30; B="코를 돼지처럼 그려놨어"; C="They just can't get my nose right!"; D=1; E=44814.77115533565; F="tangel_4" }
    @{ Row=33; A=31; B="너흰 그렇겠지 멋지게 그려 줬으니.."; C="Well, it's easy for you to say! You guys look amazing."; D=1; E=44814.77667180556; F="tangel_4" }
    @{ Row=34; A=32; B="자, 날 올려줘 위에서 잡아줄께"; C="All right, Okay, give me a boost, and I'll pull you up"; D=1; E=44814.77288965278; F="tangel_4" }
    @{ Row=35; A=33; B="가방 먼저 줘"; C="Give us the satchel first."; D=1; E=44816.00422657408; F="tangel_4" }
    @{ Row=36; A=34; B="같이 고생한 게 얼만데 아직도 날 못믿어?"; C="I can't believe that after all we've been through together, you don't trust me?"; D=1; E=44814.96343554398; F="tangel_4" }
    @{ Row=37; A=35; B="미안! 손에 짐이 많아서.."; C="Sorry, my hands are full."; D=1; E=44811.00312993056; F="tangel_4" }
    @{ Row=38; A=36; B="가방을 꼭 찾아야 한다!"; C="Retrieve that satchel at any cost!"; D=1; E=44814.96381623843; F="tangel_4" }
    @{ Row=39; A=37; B="넌 내꺼야"; C="Alone, at last"; D=1; E=44814.77643712963; F="tangel_4" }
    @{ Row=40; A=38; B="내가 옷장에 사람을 가둔거야!"; C="I've got a person in my closet!"; D=1; E=44814.77538738426; F="tangel_4" }
    @{ Row=41; A=39; B="내가 나약하다고요, 엄마?"; C="Too weak to handel myself out there, huh, mother?"; D=1; E=44814.77363428241; F="tangel_4" }
    @{ Row=42; A=40; B="깜짝 선물이 있단다!"; C="I have a big surprise"; D=1; E=44814.7587378125; F="tangel_5" }
    @{ Row=43; A=41; B="내 선물이 더 놀라울 걸?"; C="Oh I bet my surprise is bigger."; D=1; E=44814.96390042824; F="tangel_5" }
    @{ Row=44; A=42; B="오 엄마, 저도 드릴 말씀 있어요"; C="Well mother, there's something I want to tell you."; D=1; E=44816.0040205787; F="tangel_5" }
    @{ Row=45; A=43; B="엄마가 했던말 생각 해 봤는데`n[많이 생각해 봤다 전에 너가 했던말에 대하여]"; C="Okay, I've been thinking a lot about what you said, earlier."; D=1; E=44816.00758071759; F="tangel_5" }
    @{ Row=46; A=44; B="그 얘긴 관뒀으면 좋겠다"; C="Because I really thought we dropped the issue, sweetheart"; D=1; E=44814.76609010417; F="tangel_5" }
    @{ Row=47; A=45; B="넌 네 앞가림 하긴 너무 약해"; C="Oh darling, I know you're not strong enough to handle yourself out there"; D=1; E=44816.00548414352; F="tangel_5" }
    @{ Row=48; A=46; B="라푼젤, 그 얘긴 그만 하자"; C="Rapunzel, We're done talking about this"; D=1; E=44814.76664453703; F="tangel_5" }
    @{ Row=49; A=47; B="불빛 얘긴 그만 해!`n넌 이 탑을 떠날 수 없어 영원히!"; C="Enough for the lights, Rapunzel! You are not leaving this tower, ever!"; D=1; E=44816.00618868056; F="tangel_5" }
    @{ Row=50; A=48; B="''별' 보여달라는 것보단 낫잖아요"; C="I just thought it was a better idea than stars"; D=1; E=44814.75894799769; F="tangel_5" }
    @{ Row=51; A=49; B="너 혼자 있을 수 있겠니?"; C="You sure You'll be all right on your own?"; D=1; E=44816.00851835648; F="tangel_5" }
    @{ Row=52; A=50; B="내 주머니 어딨어"; C="Where is my satchel?"; D=1; E=44816.01733137731; F="tangel_6" }
    @{ Row=53; A=51; B="감춰놨지, 절대 못 찾을 곳에"; C="I've hidden it, somewhere you'll never find it."; D=1; E=44816.02475321759; F="tangel_6" }
    @{ Row=54; A=52; B="저 항아리에?"; C="it's in that pot, isnt it?"; D=1; E=44816.02421196759; F="tangel_6" }
    @{ Row=55; A=53; B="내 머릴 어쩌려는 거야? 잘라가려고?"; C="So, What do you want, with my hair, to cut it?"; D=1; E=44816.02448789352; F="tangel_6" }
    @{ Row=56; A=54; B="아니야, 내가 원하는건 이 머리카락에서 벗어나는 거야, 진짜로!"; C="No! Listen, the only thing I want to do with your hair, is to get out of it. Literally"; D=1; E=44816.02268734953; F="tangel_6" }
    @{ Row=57; A=55; B="당신 머리칼을 왜 노려? `n[왜 지구에서 내가 너의 머리칼을 원해?]"; C="Why on earth would I want your hair?"; D=1; E=44816.02372802083; F="tangel_6" }
    @{ Row=58; A=56; B="그럼 그곳으로 날 안내해 줬다가 집으로 데려다 줘"; C="Take me to these lanterns, and return me home safely."; D=1; E=44816.02175428241; F="tangel_6" }
    @{ Row=59; A=57; B="이 탑을 다 꺠부수고 박살내도"; C="You can tear this tower apart, brick by brick"; D=1; E=44816.0238262037; F="tangel_6" }
    @{ Row=60; A=58; B="내 도움 없인 그 귀한 주머니 절대 못 찾아"; C="But without my help, you will never find your precious satchel."; D=1; E=44816.02117099537; F="tangel_6" }
    @{ Row=61; A=59; B="난 한번 약속하면 절대로 어기지 않아"; C="I promise. And when I promise something, I never ever break that promise."; D=1; E=44816.02362364584; F="tangel_6" }
    @{ Row=62; A=60; B="소용없어 아무것도 안 보여"; C="It's no use. I can't see anything."; D=0; E=$null; F="tangel_7" }
    @{ Row=63; A=61; B="다 내탓이야. 엄마 말을 진작 들을걸.."; C="This is all my fault. She was right"; D=1; E=44816.03170714121; F="tangel_7" }
    @{ Row=64; A=62; B="내 본명은 유진 피츠허버트야 말해주고 싶었어"; C="My real name is Eugene Fitzherbert. Someone might as well know."; D=1; E=44816.03046299768; F="tangel_7" }
    @{ Row=65; A=63; B="난 노랠 하면 머리에서 빛이 나"; C="I have magic hair that glows when I sing."; D=1; E=44816.03134275463; F="tangel_7" }
    @{ Row=66; A=64; B="살았어 살았다고!"; C="I'm alive. I'm alive!"; D=1; E=44816.03107152778; F="tangel_7" }
    @{ Row=67; A=65; B="안믿었는데."; C="I didn't see that coming."; D=0; E=$null; F="tangel_7" }
    @{ Row=68; A=66; B="진짜로 막 빛이나`n[진짜 그녀 머리에서 빛이나]"; C="Her hair actually glows."; D=0; E=$null; F="tangel_7" }
    @{ Row=69; A=67; B="왜 빛이 나지`n[왜 그녀 머리가 빛나지??]"; C="Why does her hair glow?"; D=0; E=$null; F="tangel_7" }
    @{ Row=70; A=68; B="빛만 나는게 아니야"; C="It doesn't just glow."; D=0; E=$null; F="tangel_7" }
    @{ Row=71; A=69; B="얘 왜 날 비웃지?"; C="Why is he smiling at me?"; D=1; E=44816.03193355976; F="tangel_7" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    if ($null -ne $r.E) {
        $ws.Range("E$rowNum").Value = $r.E
    }
    $ws.Range("F$rowNum").Value = $r.F

    # Match formatting used by the existing data rows: column A is bold/
    # bordered/centered (copied from A31), column E uses the date-time number
    # format (copied from E31) but only for rows that actually have a time value.
    $ws.Range("A31").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)
    if ($null -ne $r.E) {
        $ws.Range("E31").Copy()
        $ws.Range("E$rowNum").PasteSpecial(-4122)
    }
}

Write-Host "Edit applied."
